$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.536.01"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.069.06"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'231.28"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'0.622"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'57.82"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").Value = "'0.387"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").Value = "'0.0773"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.377.34"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'14.77"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'21.15"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "'0.764"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "'5.31"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "2.067.28"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "37.540.97"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'69.90"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").Value = "'226.90"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("D26").Value = "'9.87"
$ws.Range("E26").Value = "  +6.36%  "
$ws.Range("D27").Value = "'169.41"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("D29").Value = "'19.30"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").Value = "'1.35"
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "'4.56"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("D33").Value = "'0.0625"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").Value = "'4.63"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'2.53"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "'5.30"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("D41").Value = "'98.17"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").Value = "'0.0958"
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "'2.90"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.484.44"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("D45").Value = "'1.19"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").Value = "'16.64"
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("D47").Value = "'4.01"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("D49").Value = "'7.25"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "2.261.54"
$ws.Range("E51").Value = "  -0.69%  "

$resetCells = @("D5","D6","D8","D9","D10","D13","D14","D15","D16","D19","D20","D22","D25","D26","D27","D29","D30","D32","D33","D34","D35","D39","D41","D42","D43","D45","D46","D47","D49")
foreach ($addr in $resetCells) {
    $ws.Range($addr).Style = "Normal"
}
